# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# Sheet "Hoja1" contains the daily conversion note in A1
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.6 = 5864.0 pesos`n✅ 5864.0 pesos = 1.59 = 944.07 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# Sheet "tasas" holds the rate table with manually updated figures in N10/O10/N12/O12
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 625
$wsTasas.Range("O10").Value = 3665
$wsTasas.Range("N12").Value = 3684
$wsTasas.Range("O12").Value = 593.1
